$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = -0.0461
$ws.Range("E2").Value = 0.0596
$ws.Range("F2").Value = 0.0784
$ws.Range("G2").Value = 0.3766959587790005
$ws.Range("H2").Value = 0.3766959587790005
$ws.Range("I2").Value = 0.3606019574522223
$ws.Range("J2").Value = 0.3178377756735956
$ws.Range("K2").Value = 559.46
$ws.Range("L2").Value = 0.508410501540335
$ws.Range("M2").Value = 592.2
$ws.Range("N2").Value = 0.09892422825070159
$ws.Range("O2").Value = 1.058520716405105
$ws.Range("P2").Value = 253.2
$ws.Range("Q2").Value = 0.0422958706401176
$ws.Range("R2").Value = 0.4525792728702677
$ws.Range("S2").Value = 339
$ws.Range("T2").Value = 0.5724417426545085
$ws.Range("U2").Value = 6305.047
$ws.Range("V2").Value = 1.053228484565014
$ws.Range("W2").Value = 0.1418558830012848
$ws.Range("X2").Value = 0.02064507050059136
$ws.Range("Y2").Value = 0.1212108125006935
$ws.Range("Z2").Value = -0.1733844112066068
$ws.Range("AA2").Value = -0.01029795636618661
$ws.Range("AB2").Value = 0.02079679292479834
$ws.Range("AC2").Value = -0.04364914128770801
$ws.Range("AD2").Value = 1154.72
$ws.Range("AE2").Value = 0
$ws.Range("AF2").Value = 1154.72
$ws.Range("AG2").Value = -5150.327
$ws.Range("AH2").Value = 0.1617001254705144
$ws.Range("AI2").Value = 0.2514951801403487
$ws.Range("AJ2").Value = -6.160140322675171
$ws.Range("AK2").Value = 3.005512284762086
$ws.Range("AL2").Value = 31.63
$ws.Range("AM2").Value = 31.63
$ws.Range("AN2").Value = 2.74097987086973
$ws.Range("AO2").Value = 12.54536832121404
$ws.Range("AP2").Value = -12.2254248955564
$ws.Range("AQ2").Value = 12.54536832121404
$ws.Range("D3").Value = -0.122
$ws.Range("E3").Value = 0.131
$ws.Range("K3").Value = 99.1
$ws.Range("L3").Value = 1.009164969450102
$ws.Range("M3").Value = 0
$ws.Range("N3").Value = 0
$ws.Range("O3").Value = 0
$ws.Range("P3").Value = 0
$ws.Range("Q3").Value = 0
$ws.Range("R3").Value = 0
$ws.Range("U3").Value = 0.927
$ws.Range("V3").Value = 0.002049524651779792
$ws.Range("W3").Value = 0.5930580490724117
$ws.Range("X3").Value = 0.02055145315271485
$ws.Range("Y3").Value = 0.5725065959196968
$ws.Range("Z3").Value = 0.5905310570024837
$ws.Range("AB3").Value = 0.02055145315271485
$ws.Range("AC3").Value = -0.02055145315271485
$ws.Range("AG3").Value = -0.927
$ws.Range("AJ3").Value = -0.002053733829892351
$ws.Range("AK3").Value = -0.003692153278130265
$ws.Range("B4").Value = 'Public Joint-Stock Company Investment Company IC Russ-Invest (MISX:RUSI)'
$ws.Range("D4").Value = -0.235
$ws.Range("E4").Value = -0.19
$ws.Range("G4").Value = 1.154228855721393
$ws.Range("H4").Value = 1.154228855721393
$ws.Range("I4").Value = -0.691542288557214
$ws.Range("J4").Value = -0.637960958565949
$ws.Range("K4").Value = 2.46
$ws.Range("L4").Value = 1.223880597014926
$ws.Range("M4").Value = 0
$ws.Range("N4").Value = 0
$ws.Range("O4").Value = 0
$ws.Range("P4").Value = 0
$ws.Range("Q4").Value = 0
$ws.Range("R4").Value = 0
$ws.Range("S4").Value = 0
$ws.Range("U4").Value = 1.42
$ws.Range("V4").Value = 0.02526690391459074
$ws.Range("W4").Value = 0.03710407239819005
$ws.Range("X4").Value = 0.02071481944760124
$ws.Range("Y4").Value = 0.01638925295058881
$ws.Range("Z4").Value = 0.03228397044651461
$ws.Range("AA4").Value = -0.02059591273237323
$ws.Range("AB4").Value = 0.02099089956297952
$ws.Range("AC4").Value = -0.04158681229535276
$ws.Range("AD4").Value = 1.02
$ws.Range("AF4").Value = 1.02
$ws.Range("AG4").Value = -0.3999999999999999
$ws.Range("AH4").Value = 0.01782593498776651
$ws.Range("AI4").Value = 0.01657997399219766
$ws.Range("AJ4").Value = -0.007168458781362005
$ws.Range("AK4").Value = -0.00665557404326123
$ws.Range("AL4").Value = 0.03
$ws.Range("AM4").Value = 0.03
$ws.Range("AN4").Value = -0.8360655737704918
$ws.Range("AO4").Value = -46.33333333333333
$ws.Range("AP4").Value = 0.3278688524590163
$ws.Range("AQ4").Value = -46.33333333333333
$ws.Range("B5").Value = 'Public Joint Stock Company "SAFMAR Financial investments" (MISX:SFIN)'
$ws.Range("D5").Value = 0.338
$ws.Range("E5").Value = 0.388
$ws.Range("G5").Value = 0
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 148.1
$ws.Range("L5").Value = 0.3786755305548453
$ws.Range("M5").Value = 361.7
$ws.Range("N5").Value = 0.5899526993965095
$ws.Range("O5").Value = 2.442268737339635
$ws.Range("P5").Value = 22.7
$ws.Range("Q5").Value = 0.03702495514597945
$ws.Range("R5").Value = 0.1532748143146523
$ws.Range("S5").Value = 339
$ws.Range("T5").Value = 0.9372408072988665
$ws.Range("U5").Value = 66.6
$ws.Range("V5").Value = 0.1086282824987767
$ws.Range("W5").Value = 0.1232933732933733
$ws.Range("X5").Value = 0.03729998716885198
$ws.Range("Y5").Value = 0.0859933861245213
$ws.Range("Z5").Value = 0.2789188418199972
$ws.Range("AA5").Value = 0
$ws.Range("AB5").Value = 0.04571147028006327
$ws.Range("AC5").Value = -0.04571147028006327
$ws.Range("AD5").Value = 1140.8
$ws.Range("AF5").Value = 1140.8
$ws.Range("AG5").Value = 1074.2
$ws.Range("AH5").Value = 0.6504361708193169
$ws.Range("AI5").Value = 0.4451207616372079
$ws.Range("AJ5").Value = 0.6366384164049073
$ws.Range("AK5").Value = 0.4303168689660697
$ws.Range("AM5").Value = 0
$ws.Range("D6").Value = 0.0298
$ws.Range("E6").Value = -0.0118
$ws.Range("F6").Value = 0.0784
$ws.Range("G6").Value = 0.6767361681168937
$ws.Range("H6").Value = 0.6767361681168937
$ws.Range("I6").Value = 0.6537514365457232
$ws.Range("J6").Value = 0.521789251619913
$ws.Range("K6").Value = 309.8
$ws.Range("L6").Value = 0.508619274339189
$ws.Range("M6").Value = 230.5
$ws.Range("N6").Value = 0.04738118730471962
$ws.Range("O6").Value = 0.7440284054228534
$ws.Range("P6").Value = 230.5
$ws.Range("Q6").Value = 0.04738118730471962
$ws.Range("R6").Value = 0.7440284054228534
$ws.Range("U6").Value = 6236.1
$ws.Range("V6").Value = 1.281882091761223
$ws.Range("W6").Value = 0.1604183927091964
$ws.Range("X6").Value = 0.02057532155358148
$ws.Range("Y6").Value = 0.1398430711556149
$ws.Range("Z6").Value = -0.07635319778373906
$ws.Range("AA6").Value = -0.0398402779303644
$ws.Range("AB6").Value = 0.02060268628661715
$ws.Range("AC6").Value = -0.06044296421698155
$ws.Range("AD6").Value = 12.9
$ws.Range("AE6").Value = 0
$ws.Range("AF6").Value = 12.9
$ws.Range("AG6").Value = -6223.200000000001
$ws.Range("AH6").Value = 0.002644689095270312
$ws.Range("AI6").Value = 0.007521865889212828
$ws.Range("AJ6").Value = 4.581272084805653
$ws.Range("AK6").Value = 1.376479175421911
$ws.Range("AL6").Value = 31.6
$ws.Range("AM6").Value = 31.6
$ws.Range("AN6").Value = 0.03053254437869823
$ws.Range("AO6").Value = 12.60126582278481
$ws.Range("AP6").Value = 12.60126582278481

$ws.Range("T3").ClearContents()
$ws.Range("T4").ClearContents()
$ws.Range("AN5").ClearContents()
$ws.Range("AP5").ClearContents()
$ws.Range("AQ5").ClearContents()
